$wb = $excel.ActiveWorkbook

# New B-column values to fill in for row 100 and row 101 (previously 0),
# plus a brand-new row 102 (date 45961, remn_amt 0), one set per sheet.
$newValues = @(
    @{ B100 = 457234; B101 = 441117 },
    @{ B100 = 85150;  B101 = 77034  },
    @{ B100 = 15426;  B101 = 14644  }
)

for ($i = 1; $i -le 3; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $vals = $newValues[$i - 1]

    # Update existing rows 100 and 101, column B
    $ws.Cells.Item(100, 2).Value = $vals.B100
    $ws.Cells.Item(101, 2).Value = $vals.B101

    # Add new row 102: date in column A, remn_amt 0 in column B
    $ws.Cells.Item(102, 1).Value = 45961
    $ws.Cells.Item(102, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item(102, 2).Value = 0
}
